# Commit: "Added a picture for the ad"
#
# The canonical OOXML diff for this change only shows two real content
# edits inside the deck (the remaining hunks are just namespace-attribute
# re-ordering noise on two otherwise-untouched extLst blobs -- an
# artifact of the original authoring tool's XML serializer, not of any
# PowerPoint object-model action, and there's no Guides/ThemeFamily COM
# surface to reach them from here):
#
#   1. Slide 6 title "Scenario / Placing a new auction": the runs
#      "<a:t> </a:t>" and "<a:t>a </a:t>" collapse into a single run
#      "<a:t> a </a:t>" (same run properties either side).
#   2. Slide 9 title "Sceanrio - ... wait for the confirmation message":
#      the runs "<a:t>fo</a:t>" and "<a:t>r</a:t>" collapse into a single
#      run "<a:t>for</a:t>" (same run properties either side).
#
# Both are reproduced by re-typing the exact same characters across the
# old run boundary via TextRange.Characters(start, length) -- exactly
# what PowerPoint does internally when you select text spanning a run
# boundary and retype it identically: the selected runs merge into one.

$p = $ppt.ActivePresentation

# --- Slide 6: "Scenario" + break + "Placing a new auction" -------------
$slide6 = $p.Slides.Item(6)
$title6 = $slide6.Shapes.Item(1)
$tr6 = $title6.TextFrame.TextRange
$full6 = $tr6.Text
$needle6 = " a "
$zeroIdx6 = $full6.IndexOf($needle6)
$run6 = $tr6.Characters($zeroIdx6 + 1, $needle6.Length)
$run6.Text = $needle6

# --- Slide 9: "Sceanrio - Make another bid and wait for the confirmation
#               message" -------------------------------------------------
$slide9 = $p.Slides.Item(9)
$title9 = $slide9.Shapes.Item(1)
$tr9 = $title9.TextFrame.TextRange
$full9 = $tr9.Text
$needle9 = "for"
$zeroIdx9 = $full9.IndexOf($needle9)
$run9 = $tr9.Characters($zeroIdx9 + 1, $needle9.Length)
$run9.Text = $needle9
